# Apply updates described by the diff:
# - Rename sheet "Through 2022-07-21" -> "Through 2022-07-22"
# - Update header string "2022 (through 07-21)" -> "2022 (through 07-22)"
# - Update I8 (July total) 118 -> 125
# - Update I14 (Total) 924 -> 931

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Name = "Through 2022-07-22"

$ws.Range("I1").Value = "2022 (through 07-22)"

$ws.Range("I8").Value = 125

$ws.Range("I14").Value = 931
